# Scheduled-runner market-data refresh for the Mandragora leve-profit sheets.
# For each listed (sheet, row) the currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are overwritten with freshly pulled values; a couple of cells that
# no longer have a value in the refreshed data are cleared outright.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 7370
$ws.Range("I62").Value = 7370
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7370
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6746
# Row 65
$ws.Range("H65").Value = 7370
$ws.Range("I65").Value = 7370
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 36850
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -33730
# Row 98
$ws.Range("H98").Value = 1888.9565
$ws.Range("I98").Value = 1937.5
$ws.Range("J98").Value = 1778
$ws.Range("K98").Value = 1937.5
$ws.Range("L98").Value = 1778
$ws.Range("M98").Value = -439.5
$ws.Range("N98").Value = -4774
# Row 122
$ws.Range("H122").Value = 1888.9565
$ws.Range("I122").Value = 1937.5
$ws.Range("J122").Value = 1778
$ws.Range("K122").Value = 5812.5
$ws.Range("L122").Value = 5334
$ws.Range("M122").Value = -3362.5
$ws.Range("N122").Value = -10234
# Row 137
$ws.Range("H137").Value = 1828.2903
$ws.Range("J137").Value = 1587.2941
$ws.Range("L137").Value = 4761.8823
$ws.Range("N137").Value = -9861.882300000001
$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 1500
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = 1500
$ws.Range("N19").Value = -1958
# Row 22
$ws.Range("H22").Value = 1750
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1451
# Row 26
$ws.Range("H26").Value = 890877.75
$ws.Range("I26").Value = 890877.75
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 890877.75
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -890547.75
# Row 45
$ws.Range("H45").Value = 2949
$ws.Range("I45").Value = 1771.3846
$ws.Range("J45").Value = 3905.8125
$ws.Range("K45").Value = 1771.3846
$ws.Range("L45").Value = 3905.8125
$ws.Range("M45").Value = -1394.3846
$ws.Range("N45").Value = -4659.8125
# Row 122
$ws.Range("H122").Value = 2403.55
$ws.Range("I122").Value = 1193.1875
$ws.Range("J122").Value = 7245
$ws.Range("K122").Value = 3579.5625
$ws.Range("L122").Value = 21735
$ws.Range("M122").Value = -1129.5625
$ws.Range("N122").Value = -26635
# Row 133
$ws.Range("H133").Value = 36747.5
$ws.Range("J133").Value = 36747.5
$ws.Range("L133").Value = 36747.5
$ws.Range("N133").Value = -41807.5
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1525.1936
$ws.Range("I122").Value = 804.5909
$ws.Range("J122").Value = 3286.6667
$ws.Range("K122").Value = 2413.7727
$ws.Range("L122").Value = 9860.000100000001
$ws.Range("M122").Value = 36.22730000000001
$ws.Range("N122").Value = -14760.0001
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 186.94737
$ws.Range("I12").Value = 70.2
$ws.Range("J12").Value = 228.64285
$ws.Range("K12").Value = 210.6
$ws.Range("L12").Value = 685.9285500000001
$ws.Range("M12").Value = -37.60000000000002
$ws.Range("N12").Value = -1031.92855
# Row 23
$ws.Range("H23").Value = 11355.556
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 11355.556
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 34066.66800000001
$ws.Range("N23").Value = -34536.66800000001
# Row 107
$ws.Range("H107").Value = 841713.4399999999
$ws.Range("I107").Value = 481.2
$ws.Range("K107").Value = 1443.6
$ws.Range("M107").Value = 476.4000000000001
# Row 122
$ws.Range("H122").Value = 3791.7842
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3791.7842
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 34126.0578
$ws.Range("N122").Value = -39026.0578
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 31765.857
$ws.Range("J4").Value = 71669.336
$ws.Range("L4").Value = 71669.336
$ws.Range("N4").Value = -71893.336
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2584.8572
$ws.Range("I7").Value = 2267
$ws.Range("J7").Value = 3157
$ws.Range("K7").Value = 2267
$ws.Range("L7").Value = 3157
$ws.Range("M7").Value = -2155
$ws.Range("N7").Value = -3381
# Row 20
$ws.Range("H20").Value = 26000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30452
# Row 22
$ws.Range("H22").Value = 1450.125
$ws.Range("I22").Value = 1491.1818
$ws.Range("J22").Value = 1359.8
$ws.Range("K22").Value = 1491.1818
$ws.Range("L22").Value = 1359.8
$ws.Range("M22").Value = -1196.1818
$ws.Range("N22").Value = -1949.8
# Row 27
$ws.Range("H27").Value = 1450.125
$ws.Range("I27").Value = 1491.1818
$ws.Range("J27").Value = 1359.8
$ws.Range("K27").Value = 1491.1818
$ws.Range("L27").Value = 1359.8
$ws.Range("M27").Value = -1384.1818
$ws.Range("N27").Value = -1573.8
# Row 122
$ws.Range("H122").Value = 10075.474
$ws.Range("I122").Value = 10965.692
$ws.Range("J122").Value = 8146.6665
$ws.Range("K122").Value = 32897.076
$ws.Range("L122").Value = 24439.9995
$ws.Range("M122").Value = -30447.076
$ws.Range("N122").Value = -29339.9995
# Row 126
$ws.Range("H126").Value = 2584.8572
$ws.Range("I126").Value = 2267
$ws.Range("J126").Value = 3157
$ws.Range("K126").Value = 6801
$ws.Range("L126").Value = 9471
$ws.Range("M126").Value = -4331
$ws.Range("N126").Value = -14411
$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -71
# Row 81
$ws.Range("H81").Value = 2500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = 5000
$ws.Range("N81").Value = -7122
# Row 84
$ws.Range("H84").Value = 2500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = 25000
$ws.Range("N84").Value = -35608
# Row 107
$ws.Range("H107").Value = 76923330
$ws.Range("I107").Value = 252
$ws.Range("K107").Value = 756
$ws.Range("M107").Value = 1164
# Row 122
$ws.Range("H122").Value = 5786.6294
$ws.Range("I122").Value = 6147.4585
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 18442.3755
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -15992.3755
$ws.Range("N122").Value = -13600
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0
# Row 136
$ws.Range("H136").Value = 5320016.5
$ws.Range("I136").Value = 5814701.5
$ws.Range("J136").Value = 2150
$ws.Range("K136").Value = 17444104.5
$ws.Range("L136").Value = 6450
$ws.Range("M136").Value = -17441554.5
$ws.Range("N136").Value = -11550

Write-Output "Applied 190 cell edits across $($wb.Worksheets.Count) worksheets"
